$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend column A's numeric-index styling (bold/border/centered, same as
# A2:A7) down through the newly added rows 8:15 before filling values.
$ws.Range("A7").Copy()
$ws.Range("A8:A15").PasteSpecial(-4122)

$data = @(
    @(0,  "NSE:ALKYLAMINE", "NSE:360ONE",     "NSE:BHARTIARTL", "NSE:DLF", "NSE:ASIANPAINT"),
    @(1,  "NSE:BHARTIARTL", "NSE:ALLCARGO",   "NSE:ITC",        "",        "NSE:BHARTIARTL"),
    @(2,  "NSE:DALBHARAT",  "NSE:CCHHL",      "",               "",        "NSE:HINDUNILVR"),
    @(3,  "NSE:ESCORTS",    "NSE:CELLO",      "",               "",        ""),
    @(4,  "NSE:GOLDETF",    "NSE:EICHERMOT",  "",               "",        ""),
    @(5,  "NSE:HINDUNILVR", "NSE:GMBREW",     "",               "",        ""),
    @(6,  "NSE:JYOTHYLAB",  "NSE:INDBANK",    "",               "",        ""),
    @(7,  "",               "NSE:JUBLFOOD",   "",               "",        ""),
    @(8,  "",               "NSE:KPIGREEN",   "",               "",        ""),
    @(9,  "",               "NSE:LEMONTREE",  "",               "",        ""),
    @(10, "",               "NSE:MON100",     "",               "",        ""),
    @(11, "",               "NSE:MRF",        "",               "",        ""),
    @(12, "",               "NSE:PRESTIGE",   "",               "",        ""),
    @(13, "",               "NSE:ROML",       "",               "",        "")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
}
